$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.069915883399338
$ws.Range("D2").Value = 1.069537096958471
$ws.Range("E2").Value = 1.074422159387467
$ws.Range("F2").Value = 1.085481800360076
$ws.Range("I2").Value = 1.06020519710477
$ws.Range("J2").Value = 1.074848185507594
$ws.Range("K2").Value = 1.072239339789143
$ws.Range("L2").Value = 1.077111414247328
$ws.Range("M2").Value = 1.08814211713899
$ws.Range("N2").Value = 1.076374593217824
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.071279596921514
$ws.Range("D3").Value = 1.070633905237736
$ws.Range("E3").Value = 1.075685197148978
$ws.Range("F3").Value = 1.086884468388393
$ws.Range("I3").Value = 1.060734663005076
$ws.Range("J3").Value = 1.075867286217725
$ws.Range("K3").Value = 1.073151641589703
$ws.Range("L3").Value = 1.078190479479568
$ws.Range("M3").Value = 1.08936258244767
$ws.Range("N3").Value = 1.077395141167856
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.072160606606867
$ws.Range("D4").Value = 1.071342357352474
$ws.Range("E4").Value = 1.076500885043493
$ws.Range("F4").Value = 1.087790901154161
$ws.Range("I4").Value = 1.061075260231069
$ws.Range("J4").Value = 1.076524809762601
$ws.Range("K4").Value = 1.073740112153179
$ws.Range("L4").Value = 1.07888656478402
$ws.Range("M4").Value = 1.090150581425795
$ws.Range("N4").Value = 1.078053598471586
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.072530652811515
$ws.Range("D5").Value = 1.071639893969491
$ws.Range("E5").Value = 1.076843426870658
$ws.Range("F5").Value = 1.088171687053561
$ws.Range("I5").Value = 1.061217970483624
$ws.Range("J5").Value = 1.076800781628609
$ws.Range("K5").Value = 1.073987066268621
$ws.Range("L5").Value = 1.079178691228416
$ws.Range("M5").Value = 1.090481448865218
$ws.Range("N5").Value = 1.078329962249298
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.072592766001177
$ws.Range("D6").Value = 1.071689834399045
$ws.Range("E6").Value = 1.076900919412013
$ws.Range("F6").Value = 1.088235606521867
$ws.Range("I6").Value = 1.061241904277506
$ws.Range("J6").Value = 1.076847092200607
$ws.Range("K6").Value = 1.074028505326721
$ws.Range("L6").Value = 1.079227710909887
$ws.Range("M6").Value = 1.090536979232313
$ws.Range("N6").Value = 1.078376338587619
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.072165552472773
$ws.Range("D7").Value = 1.071346334211891
$ws.Range("E7").Value = 1.076505463564543
$ws.Range("F7").Value = 1.087795990320844
$ws.Range("I7").Value = 1.061077169002799
$ws.Range("J7").Value = 1.076528499078858
$ws.Range("K7").Value = 1.07374341368592
$ws.Range("L7").Value = 1.078890470181041
$ws.Range("M7").Value = 1.090155004086635
$ws.Range("N7").Value = 1.078057293027095
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.070377051061585
$ws.Range("D8").Value = 1.06990803101624
$ws.Range("E8").Value = 1.074849339714175
$ws.Range("F8").Value = 1.085956087908839
$ws.Range("I8").Value = 1.060384548814287
$ws.Range("J8").Value = 1.075192992053997
$ws.Range("K8").Value = 1.072548041810008
$ws.Range("L8").Value = 1.07747653593518
$ws.Range("M8").Value = 1.088554939216701
$ws.Range("N8").Value = 1.076719889429074
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.067214464288518
$ws.Range("D9").Value = 1.067363752444291
$ws.Range("E9").Value = 1.071918680411072
$ws.Range("F9").Value = 1.082704591622791
$ws.Range("I9").Value = 1.059148608503097
$ws.Range("J9").Value = 1.072824883158261
$ws.Range("K9").Value = 1.070427301203061
$ws.Range("L9").Value = 1.074968373076563
$ws.Range("M9").Value = 1.085721979869472
$ws.Range("N9").Value = 1.07434841754708
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.065098292635125
$ws.Range("D10").Value = 1.065660715462945
$ws.Range("E10").Value = 1.069956259612144
$ws.Range("F10").Value = 1.080530272829771
$ws.Range("I10").Value = 1.058314091677511
$ws.Range("J10").Value = 1.071235925650689
$ws.Range("K10").Value = 1.069003582744606
$ws.Range("L10").Value = 1.07328476765121
$ws.Range("M10").Value = 1.083823967089555
$ws.Range("N10").Value = 1.072757203537591
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.064180038121616
$ws.Range("D11").Value = 1.064921598224397
$ws.Range("E11").Value = 1.069104382406981
$ws.Range("F11").Value = 1.079587108316737
$ws.Range("I11").Value = 1.057950195533472
$ws.Range("J11").Value = 1.070545401980181
$ws.Range("K11").Value = 1.068384694992796
$ws.Range("L11").Value = 1.072552950187295
$ws.Range("M11").Value = 1.082999807913795
$ws.Range("N11").Value = 1.072065699244267
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.063838659042831
$ws.Range("D12").Value = 1.064646797999262
$ws.Range("E12").Value = 1.068787629715386
$ws.Range("F12").Value = 1.079236516945995
$ws.Range("I12").Value = 1.05781464268956
$ws.Range("J12").Value = 1.070288530367568
$ws.Range("K12").Value = 1.068154445814508
$ws.Range("L12").Value = 1.072280692955706
$ws.Range("M12").Value = 1.082693325289952
$ws.Range("N12").Value = 1.0718084628445
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.063911899637622
$ws.Range("D13").Value = 1.064705755391197
$ws.Range("E13").Value = 1.068855589173383
$ws.Range("F13").Value = 1.07931173177909
$ws.Range("I13").Value = 1.057843736737437
$ws.Range("J13").Value = 1.070343647503478
$ws.Range("K13").Value = 1.06820385171832
$ws.Range("L13").Value = 1.072339112507605
$ws.Range("M13").Value = 1.082759082942985
$ws.Range("N13").Value = 1.071863658253065
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.064151825751344
$ws.Range("D14").Value = 1.064898888468891
$ws.Range("E14").Value = 1.069078206251594
$ws.Range("F14").Value = 1.079558133630919
$ws.Range("I14").Value = 1.057938998578553
$ws.Range("J14").Value = 1.070524176684593
$ws.Range("K14").Value = 1.068365670033144
$ws.Range("L14").Value = 1.072530454079311
$ws.Range("M14").Value = 1.082974481208171
$ws.Range("N14").Value = 1.072044443806325
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.064299612380412
$ws.Range("D15").Value = 1.065017849635355
$ws.Range("E15").Value = 1.069215324375141
$ws.Range("F15").Value = 1.079709915531542
$ws.Range("I15").Value = 1.057997641346897
$ws.Range("J15").Value = 1.070635356099825
$ws.Range("K15").Value = 1.06846532289768
$ws.Range("L15").Value = 1.072648289084174
$ws.Range("M15").Value = 1.08310714818664
$ws.Range("N15").Value = 1.072155781109078
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.065159192707366
$ws.Range("D16").Value = 1.065709732176412
$ws.Range("E16").Value = 1.070012750325941
$ws.Range("F16").Value = 1.080592831754705
$ws.Range("I16").Value = 1.058338188369736
$ws.Range("J16").Value = 1.071281700398028
$ws.Range("K16").Value = 1.069044605129929
$ws.Range("L16").Value = 1.073333276336179
$ws.Range("M16").Value = 1.083878614649804
$ws.Range("N16").Value = 1.072803043290322
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.065697860701244
$ws.Range("D17").Value = 1.066143275485897
$ws.Range("E17").Value = 1.070512377994875
$ws.Range("F17").Value = 1.081146209595837
$ws.Range("I17").Value = 1.058551120883643
$ws.Range("J17").Value = 1.071686463292798
$ws.Range("K17").Value = 1.069407325362979
$ws.Range("L17").Value = 1.073762195351222
$ws.Range("M17").Value = 1.084361912859393
$ws.Range("N17").Value = 1.073208380994829
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.066011870112496
$ws.Range("D18").Value = 1.066395991156345
$ws.Range("E18").Value = 1.070803596715612
$ws.Range("F18").Value = 1.081468824871858
$ws.Range("I18").Value = 1.05867507549825
$ws.Range("J18").Value = 1.071922314330506
$ws.Range("K18").Value = 1.069618661908836
$ws.Range("L18").Value = 1.074012106315318
$ws.Range("M18").Value = 1.084643590373192
$ws.Range("N18").Value = 1.073444566968064
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.066118907789567
$ws.Range("D19").Value = 1.066482133208361
$ws.Range("E19").Value = 1.07090286009278
$ws.Range("F19").Value = 1.081578801309269
$ws.Range("I19").Value = 1.058717299306982
$ws.Range("J19").Value = 1.072002692846882
$ws.Range("K19").Value = 1.069690683000686
$ws.Range("L19").Value = 1.074097273844699
$ws.Range("M19").Value = 1.084739597643299
$ws.Range("N19").Value = 1.073525059631152
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.065640086069924
$ws.Range("D20").Value = 1.066096777270065
$ws.Range("E20").Value = 1.070458793994995
$ws.Range("F20").Value = 1.081086854061507
$ws.Range("I20").Value = 1.058528300639025
$ws.Range("J20").Value = 1.071643060963107
$ws.Range("K20").Value = 1.069368432951246
$ws.Range("L20").Value = 1.073716204398657
$ws.Range("M20").Value = 1.084310082563182
$ws.Range("N20").Value = 1.073164917028852
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.064081181788653
$ws.Range("D21").Value = 1.064842022787015
$ws.Range("E21").Value = 1.069012660166758
$ws.Range("F21").Value = 1.079485581608862
$ws.Range("I21").Value = 1.05791095701195
$ws.Range("J21").Value = 1.070471025885674
$ws.Range("K21").Value = 1.068318028735177
$ws.Range("L21").Value = 1.072474120610429
$ws.Range("M21").Value = 1.082911061591233
$ws.Range("N21").Value = 1.071991217527175
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.063099306318668
$ws.Range("D22").Value = 1.064051606779446
$ws.Range("E22").Value = 1.068101519899803
$ws.Range("F22").Value = 1.078477301334055
$ws.Range("I22").Value = 1.057520575890034
$ws.Range("J22").Value = 1.069731917379229
$ws.Range("K22").Value = 1.067655473784828
$ws.Range("L22").Value = 1.071690695938081
$ws.Range("M22").Value = 1.08202939471664
$ws.Range("N22").Value = 1.071251059401874
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.063619983551865
$ws.Range("D23").Value = 1.064470765358088
$ws.Range("E23").Value = 1.068584714785103
$ws.Range("F23").Value = 1.079011954059635
$ws.Range("I23").Value = 1.057727737041629
$ws.Range("J23").Value = 1.070123943529517
$ws.Range("K23").Value = 1.068006909686559
$ws.Range("L23").Value = 1.072106241075998
$ws.Range("M23").Value = 1.082496979239327
$ws.Range("N23").Value = 1.071643642274261
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.065666192514714
$ws.Range("D24").Value = 1.066117788315022
$ws.Range("E24").Value = 1.070483006932988
$ws.Range("F24").Value = 1.08111367476844
$ws.Range("I24").Value = 1.058538612883039
$ws.Range("J24").Value = 1.071662673349784
$ws.Range("K24").Value = 1.069386007474594
$ws.Range("L24").Value = 1.073736986565014
$ws.Range("M24").Value = 1.084333503128982
$ws.Range("N24").Value = 1.073184557267367
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.068033412717013
$ws.Range("D25").Value = 1.068022697901587
$ws.Range("E25").Value = 1.072677825326871
$ws.Range("F25").Value = 1.08354632729187
$ws.Range("I25").Value = 1.059469977170662
$ws.Range("J25").Value = 1.073438876136197
$ws.Range("K25").Value = 1.07097728860442
$ws.Range("L25").Value = 1.075618796426558
$ws.Range("M25").Value = 1.086455995059185
$ws.Range("N25").Value = 1.074963282465471
